$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9

$ws.Range("M4").Value = 1.07
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 2.62

$ws.Range("V5").Value = 1.69

$ws.Range("U6").Value = 1.63

$ws.Range("V7").Value = 1.69
